# Generate Report for Handoff
#
# A new source file (be36bf8a-101b-4e9a-84b8-a0515b862c1b) has been handed
# off. This adds a new row for it to the "Overview" sheet and to each of the
# per-locale detail sheets ("zh-cn" and "de-de"), and - because it was
# handed off between the existing "684eee9c..." and "7cbb2cb4..." entries
# chronologically - it also swaps the order of those two existing rows so
# everything stays sorted by handoff time, before the final
# ".localization-config" footer row.

$wb = $excel.ActiveWorkbook

$newFileName        = "be36bf8a-101b-4e9a-84b8-a0515b862c1b.md"
$newStatus          = "Ready for handoff"
$newHandoffZh       = "be36bf8a-101b-4e9a-84b8-a0515b862c1b.5a2c5e40228ff3d7445d11d6b4c6af4d93cb32ed.zh-cn.xlf"
$newHandoffZhDate   = "2016-03-03 06:58:40"
$newHandoffDe       = "be36bf8a-101b-4e9a-84b8-a0515b862c1b.5a2c5e40228ff3d7445d11d6b4c6af4d93cb32ed.de-de.xlf"
$newHandoffDeDate   = "2016-03-03 06:58:50"
$neverDate          = "0001-01-01 00:00:00"
$includeReason      = "Include"

$mdBase   = "https://github.com/OpenLocalizationTest/oltest/blob/e2f5f6a7b8c9d0e1f2a3b4c5d6e7f8a9b0c1d2e3/e2e/$newFileName"
$xlfZhUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a1b2c3d4e5f6a7b8c9d0e1f2a3b4c5d6e7f8a9b0/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$newHandoffZh"
$xlfDeUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a1b2c3d4e5f6a7b8c9d0e1f2a3b4c5d6e7f8a9b0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$newHandoffDe"

# ---------------------------------------------------------------------------
# Sheet 1: "Overview" -- columns A (File Name), B (zh-cn), C (de-de)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

# Insert a fresh row above the ".localization-config" row (row 9) to hold
# the new entry; the old row 9 becomes row 10.
$ws1.Rows.Item(9).Insert()

# Rows 7 & 8 swap their roles: 7cbb2cb4 moves up to row 7, 684eee9c moves
# down to row 8.
$ws1.Cells.Item(7, 1).Value = "7cbb2cb4-f468-479a-965e-8eb53477f492.md"
$ws1.Cells.Item(7, 2).Value = "Ready for handoff"
$ws1.Cells.Item(7, 3).Value = "Ready for handoff"

$ws1.Cells.Item(8, 1).Value = "684eee9c-aa2c-45f0-94d6-0cc89432a192.md"
$ws1.Cells.Item(8, 2).Value = "Ready for handoff"
$ws1.Cells.Item(8, 3).Value = "Ready for handoff"

# New row 9: be36bf8a
$ws1.Cells.Item(9, 1).Value = $newFileName
$ws1.Cells.Item(9, 2).Value = $newStatus
$ws1.Cells.Item(9, 3).Value = $newStatus

# Row 10 keeps the ".localization-config" content (already shifted down by
# the Insert above); nothing else to change there.

# Rebuild all hyperlinks on this sheet in their final, correct positions
# (row-insert does not shift existing hyperlink anchors automatically).
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Cells.Item(2, 1), "https://github.com/OpenLocalizationTest/oltest/blob/58b79c846e4135b005615b366563bb010059bf16/e2e/882ab199-e563-46a5-935c-f92de8a7eb49.md", [System.Type]::Missing, [System.Type]::Missing, "882ab199-e563-46a5-935c-f92de8a7eb49.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Cells.Item(3, 1), "https://github.com/OpenLocalizationTest/oltest/blob/beb7b57e2f8af33a52a6a51a8e13cb8a2e7ad8f3/e2e/2aa438ea-f3e3-428b-aff7-543eca8ba0f6.md", [System.Type]::Missing, [System.Type]::Missing, "2aa438ea-f3e3-428b-aff7-543eca8ba0f6.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Cells.Item(4, 1), "https://github.com/OpenLocalizationTest/oltest/blob/eed213af7ca40fd417abe975e67be74d227528f8/e2e/721aa51a-e9be-4dc0-9833-32873f099577.md", [System.Type]::Missing, [System.Type]::Missing, "721aa51a-e9be-4dc0-9833-32873f099577.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Cells.Item(5, 1), "https://github.com/OpenLocalizationTest/oltest/blob/eed213af7ca40fd417abe975e67be74d227528f8/e2e/af6c4662-f8fd-4e34-957c-3654765d9d23.md", [System.Type]::Missing, [System.Type]::Missing, "af6c4662-f8fd-4e34-957c-3654765d9d23.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Cells.Item(6, 1), "https://github.com/OpenLocalizationTest/oltest/blob/7f290168dd55e74578c4fa49adf98423d1074b9f/e2e/d9f67aad-a7d6-4439-917e-cf6f80be5e91.md", [System.Type]::Missing, [System.Type]::Missing, "d9f67aad-a7d6-4439-917e-cf6f80be5e91.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Cells.Item(7, 1), "https://github.com/OpenLocalizationTest/oltest/blob/b175a3ae33acfa4e0251f8984d1e35aaf94bc735/e2e/7cbb2cb4-f468-479a-965e-8eb53477f492.md", [System.Type]::Missing, [System.Type]::Missing, "7cbb2cb4-f468-479a-965e-8eb53477f492.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Cells.Item(8, 1), "https://github.com/OpenLocalizationTest/oltest/blob/21a61a212b863ff40f65614b58cbc49a1f2fb1cc/e2e/684eee9c-aa2c-45f0-94d6-0cc89432a192.md", [System.Type]::Missing, [System.Type]::Missing, "684eee9c-aa2c-45f0-94d6-0cc89432a192.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Cells.Item(9, 1), $mdBase, [System.Type]::Missing, [System.Type]::Missing, $newFileName) | Out-Null
$ws1.Hyperlinks.Add($ws1.Cells.Item(10, 1), "https://github.com/OpenLocalizationTest/oltest/blob/58b79c846e4135b005615b366563bb010059bf16/.localization-config", [System.Type]::Missing, [System.Type]::Missing, ".localization-config") | Out-Null

# ---------------------------------------------------------------------------
# Sheets 2 & 3: "zh-cn" / "de-de" detail sheets -- same row changes,
# plus their extra "Latest Handoff File/Datetime" (C/D) and
# "Latest Handback DateTime"/"Handoff Reason" (G/H) columns.
# ---------------------------------------------------------------------------
$locales = @(
    @{ Name = "zh-cn";
       C7 = "684eee9c-aa2c-45f0-94d6-0cc89432a192.6fc5b2ed88add7ea9b90544aed0eb62b1b3acf08.zh-cn.xlf"; D7 = "2016-03-03 06:55:56";
       C8 = "7cbb2cb4-f468-479a-965e-8eb53477f492.312cc3e5efb0e0c1fd7b54275018a861b3883926.zh-cn.xlf"; D8 = "2016-03-03 06:49:09";
       C9 = $newHandoffZh; D9 = $newHandoffZhDate;
       HyperlinkC7Url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e1b9a0cc485dfc2d23a3386a82bbef09bbd71abc/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/684eee9c-aa2c-45f0-94d6-0cc89432a192.6fc5b2ed88add7ea9b90544aed0eb62b1b3acf08.zh-cn.xlf";
       HyperlinkC8Url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/57ca2768c064dfd1d12e82d540d9cdfc387c1fe4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/7cbb2cb4-f468-479a-965e-8eb53477f492.312cc3e5efb0e0c1fd7b54275018a861b3883926.zh-cn.xlf";
       HyperlinkC9Url = $xlfZhUrl;
       MdA7Url = "https://github.com/OpenLocalizationTest/oltest/blob/21a61a212b863ff40f65614b58cbc49a1f2fb1cc/e2e/684eee9c-aa2c-45f0-94d6-0cc89432a192.md";
       MdA8Url = "https://github.com/OpenLocalizationTest/oltest/blob/b175a3ae33acfa4e0251f8984d1e35aaf94bc735/e2e/7cbb2cb4-f468-479a-965e-8eb53477f492.md";
    },
    @{ Name = "de-de";
       C7 = "684eee9c-aa2c-45f0-94d6-0cc89432a192.6fc5b2ed88add7ea9b90544aed0eb62b1b3acf08.de-de.xlf"; D7 = "2016-03-03 06:56:06";
       C8 = "7cbb2cb4-f468-479a-965e-8eb53477f492.312cc3e5efb0e0c1fd7b54275018a861b3883926.de-de.xlf"; D8 = "2016-03-03 06:49:20";
       C9 = $newHandoffDe; D9 = $newHandoffDeDate;
       HyperlinkC7Url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/792fda8f1fd4c8cd7e06b42bf1ad4ddfeda729b0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/684eee9c-aa2c-45f0-94d6-0cc89432a192.6fc5b2ed88add7ea9b90544aed0eb62b1b3acf08.de-de.xlf";
       HyperlinkC8Url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3d70867b5b6d8a434b7b33dc8c73cfb2ad127abb/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/7cbb2cb4-f468-479a-965e-8eb53477f492.312cc3e5efb0e0c1fd7b54275018a861b3883926.de-de.xlf";
       HyperlinkC9Url = $xlfDeUrl;
       MdA7Url = "https://github.com/OpenLocalizationTest/oltest/blob/21a61a212b863ff40f65614b58cbc49a1f2fb1cc/e2e/684eee9c-aa2c-45f0-94d6-0cc89432a192.md";
       MdA8Url = "https://github.com/OpenLocalizationTest/oltest/blob/b175a3ae33acfa4e0251f8984d1e35aaf94bc735/e2e/7cbb2cb4-f468-479a-965e-8eb53477f492.md";
    }
)

foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale.Name)

    # Insert a fresh row above the ".localization-config" row (row 9).
    $ws.Rows.Item(9).Insert()

    # Row 7: now holds 7cbb2cb4 (was 684eee9c)
    $ws.Cells.Item(7, 1).Value = "7cbb2cb4-f468-479a-965e-8eb53477f492.md"
    $ws.Cells.Item(7, 2).Value = "Ready for handoff"
    $ws.Cells.Item(7, 3).Value = $locale.C8
    $ws.Cells.Item(7, 4).Value = $locale.D8
    $ws.Cells.Item(7, 7).Value = $neverDate
    $ws.Cells.Item(7, 8).Value = $includeReason

    # Row 8: now holds 684eee9c (was 7cbb2cb4)
    $ws.Cells.Item(8, 1).Value = "684eee9c-aa2c-45f0-94d6-0cc89432a192.md"
    $ws.Cells.Item(8, 2).Value = "Ready for handoff"
    $ws.Cells.Item(8, 3).Value = $locale.C7
    $ws.Cells.Item(8, 4).Value = $locale.D7
    $ws.Cells.Item(8, 7).Value = $neverDate
    $ws.Cells.Item(8, 8).Value = $includeReason

    # Row 9 (new): be36bf8a
    $ws.Cells.Item(9, 1).Value = $newFileName
    $ws.Cells.Item(9, 2).Value = $newStatus
    $ws.Cells.Item(9, 3).Value = $locale.C9
    $ws.Cells.Item(9, 4).Value = $locale.D9
    $ws.Cells.Item(9, 7).Value = $neverDate
    $ws.Cells.Item(9, 8).Value = $includeReason

    # Row 10 keeps the ".localization-config" content, already shifted down.

    # Rebuild all hyperlinks for this sheet.
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Cells.Item(2, 1), "https://github.com/OpenLocalizationTest/oltest/blob/58b79c846e4135b005615b366563bb010059bf16/e2e/882ab199-e563-46a5-935c-f92de8a7eb49.md", [System.Type]::Missing, [System.Type]::Missing, "882ab199-e563-46a5-935c-f92de8a7eb49.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Cells.Item(2, 3), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/91110dc133eef7e81fe018d4cce2f0c2c426843f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/882ab199-e563-46a5-935c-f92de8a7eb49.86bdb211ea3f57de819e46a39902233bc715cfc3.zh-cn.xlf", [System.Type]::Missing, [System.Type]::Missing, $ws.Cells.Item(2,3).Value2) | Out-Null
    $ws.Hyperlinks.Add($ws.Cells.Item(2, 5), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/73e8cae20885087d9c8ffd4daabe474eff3b8589/e2e/882ab199-e563-46a5-935c-f92de8a7eb49.md", [System.Type]::Missing, [System.Type]::Missing, "882ab199-e563-46a5-935c-f92de8a7eb49.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Cells.Item(2, 6), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/c72f0b4954ddfdfd23b2bbd729718e68258deff6/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/882ab199-e563-46a5-935c-f92de8a7eb49.86bdb211ea3f57de819e46a39902233bc715cfc3.zh-cn.xlf", [System.Type]::Missing, [System.Type]::Missing, $ws.Cells.Item(2,6).Value2) | Out-Null

    $ws.Hyperlinks.Add($ws.Cells.Item(3, 1), "https://github.com/OpenLocalizationTest/oltest/blob/beb7b57e2f8af33a52a6a51a8e13cb8a2e7ad8f3/e2e/2aa438ea-f3e3-428b-aff7-543eca8ba0f6.md", [System.Type]::Missing, [System.Type]::Missing, "2aa438ea-f3e3-428b-aff7-543eca8ba0f6.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Cells.Item(3, 3), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b2906c50a2f54f04353797ca69041d79cb2d0fe6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/2aa438ea-f3e3-428b-aff7-543eca8ba0f6.3493345316d0d650da8b30231ef4f293442fe2f6.zh-cn.xlf", [System.Type]::Missing, [System.Type]::Missing, $ws.Cells.Item(3,3).Value2) | Out-Null

    $ws.Hyperlinks.Add($ws.Cells.Item(4, 1), "https://github.com/OpenLocalizationTest/oltest/blob/eed213af7ca40fd417abe975e67be74d227528f8/e2e/721aa51a-e9be-4dc0-9833-32873f099577.md", [System.Type]::Missing, [System.Type]::Missing, "721aa51a-e9be-4dc0-9833-32873f099577.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Cells.Item(4, 3), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1e384b61201ce780c9dd60048116ca64bb0b41c7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/721aa51a-e9be-4dc0-9833-32873f099577.3fba059ee59d5fa5ed5fd5aa2effe57f558ed525.zh-cn.xlf", [System.Type]::Missing, [System.Type]::Missing, $ws.Cells.Item(4,3).Value2) | Out-Null

    $ws.Hyperlinks.Add($ws.Cells.Item(5, 1), "https://github.com/OpenLocalizationTest/oltest/blob/eed213af7ca40fd417abe975e67be74d227528f8/e2e/af6c4662-f8fd-4e34-957c-3654765d9d23.md", [System.Type]::Missing, [System.Type]::Missing, "af6c4662-f8fd-4e34-957c-3654765d9d23.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Cells.Item(5, 3), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1e384b61201ce780c9dd60048116ca64bb0b41c7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/af6c4662-f8fd-4e34-957c-3654765d9d23.3759a12535d2c9f4036116f9969abb4278de4a85.zh-cn.xlf", [System.Type]::Missing, [System.Type]::Missing, $ws.Cells.Item(5,3).Value2) | Out-Null

    $ws.Hyperlinks.Add($ws.Cells.Item(6, 1), "https://github.com/OpenLocalizationTest/oltest/blob/7f290168dd55e74578c4fa49adf98423d1074b9f/e2e/d9f67aad-a7d6-4439-917e-cf6f80be5e91.md", [System.Type]::Missing, [System.Type]::Missing, "d9f67aad-a7d6-4439-917e-cf6f80be5e91.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Cells.Item(6, 3), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/095b4d36a9a0941578930c30259be3881f3755ea/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/d9f67aad-a7d6-4439-917e-cf6f80be5e91.278645c454e70b05689038ce54474e972a4455a8.zh-cn.xlf", [System.Type]::Missing, [System.Type]::Missing, $ws.Cells.Item(6,3).Value2) | Out-Null
    $ws.Hyperlinks.Add($ws.Cells.Item(6, 5), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/7807bfc07ffaf62b8ef6d3aefd89d9908f00b7e5/e2e/d9f67aad-a7d6-4439-917e-cf6f80be5e91.md", [System.Type]::Missing, [System.Type]::Missing, "d9f67aad-a7d6-4439-917e-cf6f80be5e91.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Cells.Item(6, 6), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/864080ded632221dd4600add044f0bbcfaf8f117/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/d9f67aad-a7d6-4439-917e-cf6f80be5e91.278645c454e70b05689038ce54474e972a4455a8.zh-cn.xlf", [System.Type]::Missing, [System.Type]::Missing, $ws.Cells.Item(6,6).Value2) | Out-Null

    $ws.Hyperlinks.Add($ws.Cells.Item(7, 1), $locale.MdA7Url, [System.Type]::Missing, [System.Type]::Missing, "7cbb2cb4-f468-479a-965e-8eb53477f492.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Cells.Item(7, 3), $locale.HyperlinkC8Url, [System.Type]::Missing, [System.Type]::Missing, $locale.C8) | Out-Null

    $ws.Hyperlinks.Add($ws.Cells.Item(8, 1), $locale.MdA8Url, [System.Type]::Missing, [System.Type]::Missing, "684eee9c-aa2c-45f0-94d6-0cc89432a192.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Cells.Item(8, 3), $locale.HyperlinkC7Url, [System.Type]::Missing, [System.Type]::Missing, $locale.C7) | Out-Null

    $ws.Hyperlinks.Add($ws.Cells.Item(9, 1), $mdBase, [System.Type]::Missing, [System.Type]::Missing, $newFileName) | Out-Null
    $ws.Hyperlinks.Add($ws.Cells.Item(9, 3), $locale.HyperlinkC9Url, [System.Type]::Missing, [System.Type]::Missing, $locale.C9) | Out-Null

    $ws.Hyperlinks.Add($ws.Cells.Item(10, 1), "https://github.com/OpenLocalizationTest/oltest/blob/58b79c846e4135b005615b366563bb010059bf16/.localization-config", [System.Type]::Missing, [System.Type]::Missing, ".localization-config") | Out-Null
}
